$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("Team A P1", "Email A P1", "Phone A P1", "Budget A P1"),
    @("Team A P2", "Email A P2", "Phone A P2", "Budget A P2"),
    @("Team B P1", "Email B P1", "Phone B P1", "Budget B P1"),
    @("Team B P2", "Email B P2", "Phone B P2", "Budget B P2")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt 4; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}
